$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 520, shifting existing rows 520:541 down to 521:542.
$ws.Rows.Item(520).Insert()

# Populate the newly inserted row 520 with the new record.
$ws.Cells.Item(520, 1).Value = 10
$ws.Cells.Item(520, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(520, 3).Value = "La Araucanía"
$ws.Cells.Item(520, 4).Value = 44509
$ws.Cells.Item(520, 5).Value = 9
$ws.Cells.Item(520, 6).Value = "Fruta"
$ws.Cells.Item(520, 7).Value = 100109
$ws.Cells.Item(520, 8).Value = "Uva"
$ws.Cells.Item(520, 9).Value = 100109001
$ws.Cells.Item(520, 10).Value = "Uva"
$ws.Cells.Item(520, 11).Value = "Superior Seedless"
$ws.Cells.Item(520, 12).Value = "Primera"
$ws.Cells.Item(520, 13).Value = 150
$ws.Cells.Item(520, 14).Value = 33000
$ws.Cells.Item(520, 15).Value = 33000
$ws.Cells.Item(520, 16).Value = 33000
$ws.Cells.Item(520, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(520, 18).Value = "EE.UU."
$ws.Cells.Item(520, 19).Value = 4125
$ws.Cells.Item(520, 20).Value = 8
